# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy formatting from the last existing header cell (AC1, bold style)
# so the new headers match the rest of row 1.
$headerSrc = $ws.Range("AC1")

$wins = $ws.Range("AD1")
$headerSrc.Copy($wins)
$wins.Value = "Wins"

$losses = $ws.Range("AE1")
$headerSrc.Copy($losses)
$losses.Value = "Losses"

$ties = $ws.Range("AF1")
$headerSrc.Copy($ties)
$ties.Value = "Ties"

# --- Data rows 2-54: Wins=80, Losses=82, Ties=0 for every player ---
$firstDataRow = 2
$lastDataRow = 54
$rowCount = $lastDataRow - $firstDataRow + 1

$dataRange = $ws.Range("AD" + $firstDataRow + ":AF" + $lastDataRow)
$values = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i, 0] = 80
    $values[$i, 1] = 82
    $values[$i, 2] = 0
}
$dataRange.Value = $values
